$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.332.14'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '1.622.87'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '212.68'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = '0.487'
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('D10').Value = '19.04'
$ws.Range('E10').Value = '  +5.64%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.850.18'
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.598.24'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('D16').Value = '26.339.09'
$ws.Range('D17').Value = '62.57'
$ws.Range('E17').Value = '  +4.21%  '
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('D21').Value = '4.31'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  +7.62%  '
$ws.Range('D25').Value = '143.43'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').Value = '15.23'
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('E29').Value = '  +1.88%  '
$ws.Range('E30').Value = '  +11.05%  '
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').Value = '  +2.95%  '
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('E34').Value = '  +3.44%  '
$ws.Range('E35').Value = '  +2.04%  '
$ws.Range('D36').Value = '1.179.45'
$ws.Range('E36').Value = '  +4.83%  '
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('D38').Value = '0.813'
$ws.Range('E38').Value = '  +3.51%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D41').Value = '0.498'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('D43').Value = '5.35'
$ws.Range('E43').Value = '  +4.65%  '
$ws.Range('D44').Value = '1.761.15'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').Value = '93.55'
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('E46').Value = '  +15.77%  '
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('D48').Value = '54.10'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').Value = '0.0509'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  -0.31%  '
